$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A12").Value = "2025/12/03 03:00"
$ws.Range("B12").Value = "-"
$ws.Range("C12").Value = "-"
$ws.Range("D12").Value = "-"
$ws.Range("E12").Value = "-"
$ws.Range("F12").Value = "-"
$ws.Range("G12").Value = "-"
